$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.526564516268195
$ws.Range("C2").Value = 0.1676714244238156
$ws.Range("D2").Value = 0.5182808915800479
$ws.Range("E2").Value = 0.1734196780717419
$ws.Range("G2").Value = 0.9526172156701591
$ws.Range("H2").Value = 0.9864382198177566
$ws.Range("J2").Value = 0.07077265090906648
$ws.Range("L2").Value = 0.4169754934032568
$ws.Range("N2").Value = 1.421750933881157
$ws.Range("O2").Value = 3.9223879105366
$ws.Range("B3").Value = 1.431762634315874
$ws.Range("C3").Value = 0.1590470105446826
$ws.Range("D3").Value = 0.5171009416604306
$ws.Range("E3").Value = 0.1742528135663459
$ws.Range("G3").Value = 0.952958548875543
$ws.Range("H3").Value = 0.9912435071558008
$ws.Range("J3").Value = 0.07114662343867773
$ws.Range("L3").Value = 0.4091006349206907
$ws.Range("N3").Value = 1.433289981010688
$ws.Range("O3").Value = 3.932836416365603
$ws.Range("B4").Value = 1.373918949326367
$ws.Range("C4").Value = 0.1537180640582108
$ws.Range("D4").Value = 0.5166032018895663
$ws.Range("E4").Value = 0.1748282724769972
$ws.Range("G4").Value = 0.9537669553418624
$ws.Range("H4").Value = 0.9946357942068005
$ws.Range("J4").Value = 0.07138888702301438
$ws.Range("L4").Value = 0.4044247209253058
$ws.Range("N4").Value = 1.440903923322935
$ws.Range("O4").Value = 3.941434158793527
$ws.Range("B5").Value = 1.350440632727896
$ws.Range("C5").Value = 0.1515381975263921
$ws.Range("D5").Value = 0.5164574850527828
$ws.Range("E5").Value = 0.1750788693155041
$ws.Range("G5").Value = 0.954246878852075
$ws.Range("H5").Value = 0.9961293431878318
$ws.Range("J5").Value = 0.07149079710288353
$ws.Range("L5").Value = 0.4025594805926289
$ws.Range("N5").Value = 1.444139792170382
$ws.Range("O5").Value = 3.94548651011209
$ws.Range("B6").Value = 1.346547769480992
$ws.Range("C6").Value = 0.1511757368767519
$ws.Range("D6").Value = 0.5164367414156459
$ws.Range("E6").Value = 0.1751214533313785
$ws.Range("G6").Value = 0.9543356565143313
$ws.Range("H6").Value = 0.9963840618414821
$ws.Range("J6").Value = 0.07150791179468241
$ws.Range("L6").Value = 0.4022521941144817
$ws.Range("N6").Value = 1.444685149915372
$ws.Range("O6").Value = 3.946192538498792
$ws.Range("B7").Value = 1.373601931944108
$ws.Range("C7").Value = 0.1536886989299546
$ws.Range("D7").Value = 0.5166010053114576
$ws.Range("E7").Value = 0.1748315869257429
$ws.Range("G7").Value = 0.9537728185428733
$ws.Range("H7").Value = 0.994655486545355
$ws.Range("J7").Value = 0.07139024851039588
$ws.Range("L7").Value = 0.4043994025021078
$ws.Range("N7").Value = 1.440947024195367
$ws.Range("O7").Value = 3.941486588633865
$ws.Range("B8").Value = 1.493802023490105
$ws.Range("C8").Value = 0.1647047722733106
$ws.Range("D8").Value = 0.5178270385706156
$ws.Range("E8").Value = 0.1736936916133445
$ws.Range("G8").Value = 0.9526105648944565
$ws.Range("H8").Value = 0.9880034208146213
$ws.Range("J8").Value = 0.07089897654626265
$ws.Range("L8").Value = 0.4142272881197613
$ws.Range("N8").Value = 1.425619893485461
$ws.Range("O8").Value = 3.925537570869977
$ws.Range("B9").Value = 1.732350142984444
$ws.Range("C9").Value = 0.1860355920552479
$ws.Range("D9").Value = 0.5220269039135275
$ws.Range("E9").Value = 0.1719684722365322
$ws.Range("G9").Value = 0.9550881539156677
$ws.Range("H9").Value = 0.9784620419844288
$ws.Range("J9").Value = 0.07003560747572335
$ws.Range("L9").Value = 0.4347569931321544
$ws.Range("N9").Value = 1.399754818978423
$ws.Range("O9").Value = 3.911584096993749
$ws.Range("B10").Value = 1.909277540573044
$ws.Range("C10").Value = 0.2015355203791103
$ws.Range("D10").Value = 0.5262034074984996
$ws.Range("E10").Value = 0.1710083930317374
$ws.Range("G10").Value = 0.9598176434188446
$ws.Range("H10").Value = 0.973585282243306
$ws.Range("J10").Value = 0.06946183216272894
$ws.Range("L10").Value = 0.4505997666961008
$ws.Range("N10").Value = 1.3832998370886
$ws.Range("O10").Value = 3.9119080427227
$ws.Range("B11").Value = 1.990115995911367
$ws.Range("C11").Value = 0.20854830313462
$ws.Range("D11").Value = 0.5283395424407331
$ws.Range("E11").Value = 0.1706381500993572
$ws.Range("G11").Value = 0.9626030556588745
$ws.Range("H11").Value = 0.9718294995154224
$ws.Range("J11").Value = 0.06921386161161402
$ws.Range("L11").Value = 0.4579706685073575
$ws.Range("N11").Value = 1.376365934843065
$ws.Range("O11").Value = 3.91435551743848
$ws.Range("B12").Value = 2.020776790017749
$ws.Range("C12").Value = 0.2111982315951764
$ws.Range("D12").Value = 0.5291823289636426
$ws.Range("E12").Value = 0.1705074916226579
$ws.Range("G12").Value = 0.9637491160463441
$ws.Range("H12").Value = 0.9712311183321276
$ws.Range("J12").Value = 0.06912183043336206
$ws.Range("L12").Value = 0.4607852590332726
$ws.Range("N12").Value = 1.37381946449689
$ws.Range("O12").Value = 3.915613270618962
$ws.Range("B13").Value = 2.01417128725825
$ws.Range("C13").Value = 0.2106277764503659
$ws.Range("D13").Value = 0.5289993144791083
$ws.Range("E13").Value = 0.1705352070620663
$ws.Range("G13").Value = 0.9634982295192032
$ws.Range("H13").Value = 0.9713570336008814
$ws.Range("J13").Value = 0.06914156792005866
$ws.Range("L13").Value = 0.4601780492720309
$ws.Range("N13").Value = 1.374364368849271
$ws.Range("O13").Value = 3.915327668886704
$ws.Range("B14").Value = 1.992637505520804
$ws.Range("C14").Value = 0.2087664287271593
$ws.Range("D14").Value = 0.5284082005325246
$ws.Range("E14").Value = 0.1706272095672468
$ws.Range("G14").Value = 0.9626955125318659
$ws.Range("H14").Value = 0.9717789378751149
$ws.Range("J14").Value = 0.06920625270805392
$ws.Range("L14").Value = 0.4582017588525815
$ws.Range("N14").Value = 1.376154847300171
$ws.Range("O14").Value = 3.91445235975354
$ws.Range("B15").Value = 1.979453774300453
$ws.Range("C15").Value = 0.2076255574819186
$ws.Range("D15").Value = 0.5280505354965186
$ws.Range("E15").Value = 0.1706848061873032
$ws.Range("G15").Value = 0.9622157169255701
$ws.Range("H15").Value = 0.9720460250246674
$ws.Range("J15").Value = 0.06924611735662634
$ws.Range("L15").Value = 0.4569942644859992
$ws.Range("N15").Value = 1.377261885641843
$ws.Range("O15").Value = 3.913959312553089
$ws.Range("B16").Value = 1.904001514528034
$ws.Range("C16").Value = 0.2010764354616015
$ws.Range("D16").Value = 0.526068551767608
$ws.Range("E16").Value = 0.1710339258911553
$ws.Range("G16").Value = 0.9596483774367499
$ws.Range("H16").Value = 0.9737093338850542
$ws.Range("J16").Value = 0.06947829959080698
$ws.Range("L16").Value = 0.4501213440860283
$ws.Range("N16").Value = 1.383764076910182
$ws.Range("O16").Value = 3.911794395231624
$ws.Range("B17").Value = 1.857803218078914
$ws.Range("C17").Value = 0.1970488563465551
$ws.Range("D17").Value = 0.5249131071841475
$ws.Range("E17").Value = 0.1712651191306627
$ws.Range("G17").Value = 0.9582358550741361
$ws.Range("H17").Value = 0.974848196784464
$ws.Range("J17").Value = 0.06962407239600488
$ws.Range("L17").Value = 0.4459468884101199
$ws.Range("N17").Value = 1.387894178795392
$ws.Range("O17").Value = 3.911055540530612
$ws.Range("B18").Value = 1.831264539648998
$ws.Range("C18").Value = 0.1947287124777688
$ws.Range("D18").Value = 0.5242707657738634
$ws.Range("E18").Value = 0.1714043559373408
$ws.Range("G18").Value = 0.9574830719255942
$ws.Range("H18").Value = 0.975546792391242
$ws.Range("J18").Value = 0.0697091450335332
$ws.Range("L18").Value = 0.4435612942071003
$ws.Range("N18").Value = 1.39032162234156
$ws.Range("O18").Value = 3.910847056468896
$ws.Range("B19").Value = 1.822284782946838
$ws.Range("C19").Value = 0.1939425406236808
$ws.Range("D19").Value = 0.5240571024503424
$ws.Range("E19").Value = 0.1714525749488018
$ws.Range("G19").Value = 0.9572384357397112
$ws.Range("H19").Value = 0.9757908060652909
$ws.Range("J19").Value = 0.06973816024564394
$ws.Range("L19").Value = 0.4427562308718791
$ws.Range("N19").Value = 1.391152432062526
$ws.Range("O19").Value = 3.910813643038779
$ws.Range("B20").Value = 1.862717665878279
$ws.Range("C20").Value = 0.197477971445835
$ws.Range("D20").Value = 0.5250338051542798
$ws.Range("E20").Value = 0.1712398604039009
$ws.Range("G20").Value = 0.9583800447535253
$ws.Range("H20").Value = 0.9747224557065124
$ws.Range("J20").Value = 0.0696084275821045
$ws.Range("L20").Value = 0.4463896694667255
$ws.Range("N20").Value = 1.387449149543066
$ws.Range("O20").Value = 3.911111786008462
$ws.Range("B21").Value = 1.99896118281498
$ws.Range("C21").Value = 0.209313306603832
$ws.Range("D21").Value = 0.5285809062886813
$ws.Range("E21").Value = 0.1705999273128995
$ws.Range("G21").Value = 0.9629288115223176
$ws.Range("H21").Value = 0.971653209983657
$ws.Range("J21").Value = 0.06918720251425903
$ws.Range("L21").Value = 0.4587816099305542
$ws.Range("N21").Value = 1.375626790100576
$ws.Range("O21").Value = 3.914700475890925
$ws.Range("B22").Value = 2.088289132947125
$ws.Range("C22").Value = 0.2170153141899505
$ws.Range("D22").Value = 0.5310965705560733
$ws.Range("E22").Value = 0.1702373181377084
$ws.Range("G22").Value = 0.9664338312207121
$ws.Range("H22").Value = 0.9700348510134233
$ws.Range("J22").Value = 0.0689228041411809
$ws.Range("L22").Value = 0.467016708434457
$ws.Range("N22").Value = 1.368362085972194
$ws.Range("O22").Value = 3.918974998278117
$ws.Range("B23").Value = 2.040587639565899
$ws.Range("C23").Value = 0.212907687989059
$ws.Range("D23").Value = 0.5297358781489692
$ws.Range("E23").Value = 0.1704257659291315
$ws.Range("G23").Value = 0.9645144045105951
$ws.Range("H23").Value = 0.9708631490710218
$ws.Range("J23").Value = 0.06906292334276642
$ws.Range("L23").Value = 0.4626090756388948
$ws.Range("N23").Value = 1.372197152267248
$ws.Range("O23").Value = 3.916517027696955
$ws.Range("B24").Value = 1.860495777425626
$ws.Range("C24").Value = 0.1972839829344366
$ws.Range("D24").Value = 0.5249791692582875
$ws.Range("E24").Value = 0.1712512601846043
$ws.Range("G24").Value = 0.9583146719085818
$ws.Range("H24").Value = 0.9747791666579673
$ws.Range("J24").Value = 0.06961549665873701
$ws.Range("L24").Value = 0.4461894434210052
$ws.Range("N24").Value = 1.387650182190626
$ws.Range("O24").Value = 3.911085683694949
$ws.Range("B25").Value = 1.667519447207439
$ws.Range("C25").Value = 0.180294767960973
$ws.Range("D25").Value = 0.5206988885975363
$ws.Range("E25").Value = 0.1723811176733321
$ws.Range("G25").Value = 0.9539076878071597
$ws.Range("H25").Value = 0.9806684372157974
$ws.Range("J25").Value = 0.07025850907740905
$ws.Range("L25").Value = 0.4290692766536637
$ws.Range("N25").Value = 1.406304043101905
$ws.Range("O25").Value = 3.913502808459612
